$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 83 (pushes the existing rows 83..182 down to 84..183)
$ws.Rows("83:83").Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(83, 1).Value  = 10
$ws.Cells.Item(83, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value  = "La Araucanía"
$ws.Cells.Item(83, 4).Value  = 44467
$ws.Cells.Item(83, 5).Value  = 9
$ws.Cells.Item(83, 6).Value  = 100112001
$ws.Cells.Item(83, 7).Value  = "Berenjena"
$ws.Cells.Item(83, 8).Value  = "Sin especificar"
$ws.Cells.Item(83, 9).Value  = "Primera"
$ws.Cells.Item(83, 10).Value = 40
$ws.Cells.Item(83, 11).Value = 12000
$ws.Cells.Item(83, 12).Value = 12000
$ws.Cells.Item(83, 13).Value = 12000
$ws.Cells.Item(83, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(83, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(83, 16).Value = 200
$ws.Cells.Item(83, 17).Value = 60
$ws.Cells.Item(83, 18).Value = "Hortaliza"
